# Apply the weekly report update:
#  - Update the "Report Generated On" timestamp
#  - Zero out the Total Billed Amount and all line-item pricing values (no-violation / re-run scenario)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# Zero out Total Billed Amount summary cell
$ws.Range("C8").Value = 0

# Zero out each line item's pricing column (H16:H23)
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0

# Zero out the TOTAL row
$ws.Range("H24").Value = 0
